# feat: add 2022-Q3 data
#
# The workbook has a "总计" (totals) sheet followed by one sheet per
# quarter (newest first): 2022-Q1, 2021-Q4, 2021-Q3.
# This change adds a new "2022-Q3" quarter:
#   - a new "2022-Q3" detail sheet is inserted right after "总计"
#     (before the existing "2022-Q1" sheet), carrying new fund data;
#   - the "总计" summary sheet gets a new row for "2022-Q3" inserted
#     above the existing "2022-Q1" row, and all following rows shift
#     down by one (their own data is unchanged).
# The existing "2022-Q1", "2021-Q4" and "2021-Q3" detail sheets and
# their data are left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q1" detail sheet, inserting the copy
#    immediately before it, then rename the copy to "2022-Q3" and
#    overwrite its data with the new quarter's figures. The untouched
#    original keeps the name "2022-Q1" and its original data.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)

$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

$q3.Range("C2").Value = "诺安全球收益不动产（QDII）"
$q3.Range("D2").Value = "'0.23"
$q3.Range("E2").Value = "'73.76"
$q3.Range("F2").Value = "'3.92"
$q3.Range("G2").Value = "'0.0090"

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row above the old
#    "2022-Q1" row (the first data row) and fill it in with the new
#    "2022-Q3" entry. The rows below shift down automatically and keep
#    their existing values; only re-number the index column (A).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Copy the index-column formatting from the (now shifted) old first
# data row onto the newly inserted row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Row-insert carries the row-above's formatting across the whole row;
# the data columns (B:D) are unstyled in every other data row, so
# strip the inherited formatting back off them.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
